# Applies the 2025-10-29 18:33 data refresh to market_health_data.xlsx
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 06:33 PM"

# --- Top Gainers sheet: rows 61-76 (Stock, Latest, Weekly, Monthly) ---
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$gainersData = @(
    @(61, "APARINDS",   3.8924, 8.3414,              15.5876),
    @(62, "HITECHGEAR", 3.8587, 1.1486,              9.9254),
    @(63, "ORIENTTECH", 3.827,  0.5247000000000001,  32.6784),
    @(64, "ICRA",       3.7985, 4.4793,              2.8828),
    @(65, "SALASAR",    3.7935, 4.7872,              11.0485),
    @(66, "NPST",       3.7841, -2.0689,             -3.5677),
    @(67, "DCW",        3.7544, 2.3219,              -3.9753),
    @(68, "RHETAN",     3.754,  4.178,               6.549),
    @(69, "HINDPETRO",  3.6935, 6.9335,              5.7397),
    @(70, "BHARTIHEXA", 3.6718, 7.0877,              15.3332),
    @(71, "HLEGLAS",    3.659,  8.115500000000001,   27.1239),
    @(72, "RHIM",       3.6544, 3.2276,              5.1826),
    @(73, "SHK",        3.6347, 2.388,               -1.932),
    @(74, "BCLIND",     3.6271, 2.2945,               0.1728),
    @(75, "MUKANDLTD",  3.6133, 11.9685,              9.550800000000001),
    @(76, "CGPOWER",    3.6125, 3.4192,               1.0325)
)
foreach ($row in $gainersData) {
    $r = $row[0]
    $wsGainers.Cells.Item($r, 2).Value = $row[1]
    $wsGainers.Cells.Item($r, 3).Value = $row[2]
    $wsGainers.Cells.Item($r, 4).Value = $row[3]
    $wsGainers.Cells.Item($r, 5).Value = $row[4]
}

# --- Top Losers sheet ---
$wsLosers = $wb.Worksheets.Item("Top Losers")

# Row 56: the Weekly column (D) now has a real figure instead of "N/A"
$wsLosers.Cells.Item(56, 4).Value = 5.2953

# Rows 70-76 (Stock, Latest, Weekly, Monthly)
$losersData = @(
    @(70, "JNKINDIA",  -2.3482, -2.8371, 4.2622),
    @(71, "FCL",       -2.3453, -2.616,  -0.02),
    @(72, "DEEDEV",    -2.3334, -6.6528, -7.4227),
    @(73, "WEALTH",    -2.2793, -3.8356, -2.7981),
    @(74, "RATNAMANI", -2.2788, -0.4626, 0.8712),
    @(75, "CSBBANK",   -2.2695, 2.3137,  10.6999),
    @(76, "BBOX",      -2.2639, -4.7636, 5.1528)
)
foreach ($row in $losersData) {
    $r = $row[0]
    $wsLosers.Cells.Item($r, 2).Value = $row[1]
    $wsLosers.Cells.Item($r, 3).Value = $row[2]
    $wsLosers.Cells.Item($r, 4).Value = $row[3]
    $wsLosers.Cells.Item($r, 5).Value = $row[4]
}

# --- 1 Month Performance sheet ---
$wsMonth = $wb.Worksheets.Item("1 Month Performance")
$wsMonth.Cells.Item(3, 3).Value = 79.9766

$monthData = @(
    @(16, "SEJALLTD", 37.4301),
    @(17, "V2RETAIL", 37.2004),
    @(18, "RAMAPHO",  36.9731),
    @(19, "SANDUMA",  36.9057)
)
foreach ($row in $monthData) {
    $r = $row[0]
    $wsMonth.Cells.Item($r, 2).Value = $row[1]
    $wsMonth.Cells.Item($r, 3).Value = $row[2]
}

$wsMonth.Cells.Item(36, 3).Value = 27.3801
